# Generate Report for Handback
# Updates timestamp / status values that are regenerated whenever the
# handback status report is (re)generated. Several cells across the
# three sheets happen to share the same text value (e.g. matching
# timestamps / statuses for the "531f0908" and "64848472" entries), so
# every cell that currently carries the old text is updated together.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- "2016-08-21 04:14:24" -> "2016-08-21 04:15:19" ---
# Overview "Latest HO Xliff Generate Date" (G) and de-de "Correspond
# Handoff Datetime" (H), rows 2 and 5.
$wsOverview.Range("G2").Value = "2016-08-21 04:15:19"
$wsOverview.Range("G5").Value = "2016-08-21 04:15:19"
$wsDeDe.Range("H2").Value = "2016-08-21 04:15:19"
$wsDeDe.Range("H5").Value = "2016-08-21 04:15:19"

# --- "ht" -> "mt" ---
# zh-cn and de-de "Priority" (E), rows 2 and 5.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# --- "2016-08-21 04:14:20" -> "2016-08-21 04:15:14" ---
# zh-cn "Correspond Handoff Datetime" (H), rows 2 and 5.
$wsZhCn.Range("H2").Value = "2016-08-21 04:15:14"
$wsZhCn.Range("H5").Value = "2016-08-21 04:15:14"

# --- "2016-08-21 04:14:44" -> "2016-08-21 04:15:32" ---
# zh-cn "Correspond Handback DateTime" (K), rows 2 and 5.
$wsZhCn.Range("K2").Value = "2016-08-21 04:15:32"
$wsZhCn.Range("K5").Value = "2016-08-21 04:15:32"

# --- "2016-08-21 04:14:50" -> "2016-08-21 04:15:38" ---
# de-de "Correspond Handback DateTime" (K), rows 2 and 5.
$wsDeDe.Range("K2").Value = "2016-08-21 04:15:38"
$wsDeDe.Range("K5").Value = "2016-08-21 04:15:38"
